$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D9").Value = "1. Execute test CAM_MOUSE_1`n2. Right click on ARE GUI background panel`n3. Click on ''Stop Model' button"
$ws.Range("D10").Value = "1. Execute test CAM_MOUSE_3`n2. Right click on ARE GUI background panel`n3. Click on 'Start Model' button"
$ws.Range("B11").Value = "Pause model/Start model"
$ws.Range("D11").Value = "1. Execute test CAM_MOUSE_1`n2. Right click on ARE GUI background panel`n3. Click on 'Pause Model' button`n4. Click on 'Start Model' button"
$ws.Range("D12").Value = "1. Execute test CAM_MOUSE_1`n2. Right click on ARE GUI background panel`n3. Click on 'Pause Model' button`n4. Click on ''Stop Model' button"

$ws.Range("D13").Select()
